$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text in H1 (was "Empresa Id", now "Entidad Id")
$ws.Range("H1").Value = "Entidad Id"

# Move the active selection to G2
$ws.Range("G2").Select()
